$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vaccinatie")

# Extend the note text in K3
$ws.Range("K3").Value2 = "Need to make one valueset containing al productcodesystems or sliced vaccineCode.coding?"

# Move the long note text from K5 to K6, clearing K5
$note = $ws.Range("K5").Value2
$ws.Range("K6").Value2 = $note
$ws.Range("K5").Value2 = ""

# Adjust row heights to match the new wrapped content
$ws.Rows.Item(3).RowHeight = 25.5
$ws.Rows.Item(5).RowHeight = 12.75
$ws.Rows.Item(6).RowHeight = 25.5

# Update the active selection on the sheet
$ws.Range("K3").Select()
